$wb = $excel.ActiveWorkbook

$logged = $wb.Worksheets.Item("Logged")
$actual = $wb.Worksheets.Item("Actual")

# --- Refresh the "Logged" sheet with freshly-scraped data -----------------
# (re-running the scraper pulled slightly different brand/byline text for
# some rows than the previous run - update those cells here)
$logged.Range("B2").Value = "Visit the Outward Hound Store"
$logged.Range("B5").Value = "Outward Hound"
$logged.Range("B6").Value = "Visit the Outward Hound Store"

# --- Compare "Logged" vs "Actual" and highlight mismatches red ------------
# Walk the shared A:E range on both sheets; any cell whose value doesn't
# match the corresponding cell on "Actual" gets filled red, matching cells
# get an (effectively invisible) cleared fill.
$lastRow = 6
$lastCol = 5

for ($r = 2; $r -le $lastRow; $r++) {
    for ($c = 1; $c -le $lastCol; $c++) {
        $loggedCell = $logged.Cells.Item($r, $c)
        $actualCell = $actual.Cells.Item($r, $c)

        if ($loggedCell.Value -eq $actualCell.Value) {
            # Matches "Actual" - clear highlight
            $loggedCell.Interior.PatternColor = 16777215
            $loggedCell.Interior.Color = 16777215
        } else {
            # Doesn't match "Actual" - flag it red
            $loggedCell.Interior.PatternColor = 15597329
            $loggedCell.Interior.Color = 15597329
        }
    }
}

# --- Restore the cursor positions left by the last interactive session ----
$logged.Range("A22").Select()
$actual.Activate()
$actual.Range("C4").Select()

Write-Host "Logged sheet cross-checked against Actual; mismatches highlighted."
